# Skill.xlsx - "modified config of skill"
# The AtkDis column (H) for every skill row (2-9) is updated from 3 to 2.5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the AtkDis values for all skill rows.
$ws.Range("H2:H9").Value = 2.5

# Rows 4-9 previously carried an explicit bordered style on column H;
# after the edit those cells fall back to the sheet's default (unstyled) look.
$ws.Range("H4:H9").Borders.LineStyle = -4142

# Leave the selection on the last edited cell, matching the authored edit.
$ws.Range("H9").Select() | Out-Null
